$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-12-25 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-12-26 Thursday", 2) | Out-Null
$d.Content.Find.Execute("34×13=442", $true, $false, $false, $false, $false, $true, 1, $false, "87×62=5394", 2) | Out-Null
$d.Content.Find.Execute("57×15=855", $true, $false, $false, $false, $false, $true, 1, $false, "25×71=1775", 2) | Out-Null
$d.Content.Find.Execute("51×30=1530", $true, $false, $false, $false, $false, $true, 1, $false, "65×84=5460", 2) | Out-Null
$d.Content.Find.Execute("17×62=1054", $true, $false, $false, $false, $false, $true, 1, $false, "26×38=988", 2) | Out-Null
$d.Content.Find.Execute("29×18=522", $true, $false, $false, $false, $false, $true, 1, $false, "69×57=3933", 2) | Out-Null
$d.Content.Find.Execute("36×32=1152", $true, $false, $false, $false, $false, $true, 1, $false, "31×56=1736", 2) | Out-Null
$d.Content.Find.Execute("21×69=1449", $true, $false, $false, $false, $false, $true, 1, $false, "39×44=1716", 2) | Out-Null
$d.Content.Find.Execute("55×70=3850", $true, $false, $false, $false, $false, $true, 1, $false, "40×97=3880", 2) | Out-Null
$d.Content.Find.Execute("51×42=2142", $true, $false, $false, $false, $false, $true, 1, $false, "67×58=3886", 2) | Out-Null
$d.Content.Find.Execute("27×25=675", $true, $false, $false, $false, $false, $true, 1, $false, "24×75=1800", 2) | Out-Null
$d.Content.Find.Execute("49×84=4116", $true, $false, $false, $false, $false, $true, 1, $false, "19×33=627", 2) | Out-Null
$d.Content.Find.Execute("18×95=1710", $true, $false, $false, $false, $false, $true, 1, $false, "65×85=5525", 2) | Out-Null
$d.Content.Find.Execute("21×79=1659", $true, $false, $false, $false, $false, $true, 1, $false, "30×71=2130", 2) | Out-Null
$d.Content.Find.Execute("74×92=6808", $true, $false, $false, $false, $false, $true, 1, $false, "86×90=7740", 2) | Out-Null
$d.Content.Find.Execute("38×61=2318", $true, $false, $false, $false, $false, $true, 1, $false, "90×96=8640", 2) | Out-Null
$d.Content.Find.Execute("60×53=3180", $true, $false, $false, $false, $false, $true, 1, $false, "34×36=1224", 2) | Out-Null
$d.Content.Find.Execute("77×23=1771", $true, $false, $false, $false, $false, $true, 1, $false, "50×61=3050", 2) | Out-Null
$d.Content.Find.Execute("24×47=1128", $true, $false, $false, $false, $false, $true, 1, $false, "21×98=2058", 2) | Out-Null
$d.Content.Find.Execute("24×37=888", $true, $false, $false, $false, $false, $true, 1, $false, "73×89=6497", 2) | Out-Null
$d.Content.Find.Execute("56×82=4592", $true, $false, $false, $false, $false, $true, 1, $false, "15×23=345", 2) | Out-Null
$d.Content.Find.Execute("37×11=407", $true, $false, $false, $false, $false, $true, 1, $false, "67×52=3484", 2) | Out-Null
$d.Content.Find.Execute("81×77=6237", $true, $false, $false, $false, $false, $true, 1, $false, "42×92=3864", 2) | Out-Null
$d.Content.Find.Execute("24×54=1296", $true, $false, $false, $false, $false, $true, 1, $false, "75×38=2850", 2) | Out-Null
$d.Content.Find.Execute("66×66=4356", $true, $false, $false, $false, $false, $true, 1, $false, "56×95=5320", 2) | Out-Null
$d.Content.Find.Execute("77×93=7161", $true, $false, $false, $false, $false, $true, 1, $false, "30×92=2760", 2) | Out-Null
